# Update countries & provincias Spain
# Applies the daily COVID data refresh: updated timestamp, refreshed case
# counts for several countries, and a handful of rows whose country
# label/order changed (Tanzania/El Salvador swap; a cluster of small
# countries around Congo/Chad/Liberia/.../Nepal/Sudan got re-sorted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Datos actualizados a ..." timestamp
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 14:42"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 85762
$ws.Range("C4").Value = 327
$ws.Range("E4").Value = 82588
$ws.Range("G4").Value = 11
$ws.Range("H4").Value = 1306

# Suecia (row 22)
$ws.Range("F22").Value = 214

# Chequia (row 26)
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 2042

# Finlandia (row 39)
$ws.Range("B39").Value = 1041
$ws.Range("C39").Value = 83
$ws.Range("E39").Value = 1026

# India (row 44)
$ws.Range("B44").Value = 863
$ws.Range("C44").Value = 136
$ws.Range("E44").Value = 770

# Principado de Andorra (row 74)
$ws.Range("F74").Value = 11

# Rows 147/148 swap: Tanzania <-> El Salvador
$ws.Range("A147").Value = "El Salvador"
$ws.Range("A148").Value = "Tanzania"
$ws.Range("D148").Value = 1
$ws.Range("E148").Value = 12

# Rows 183-192 re-sorted cluster
$ws.Range("A183").Value = "Nepal"
$ws.Range("B183").Value = 4
$ws.Range("C183").Value = 1
$ws.Range("D183").Value = 1

$ws.Range("A187").Value = "San Martin (Parte Holandesa)"

$ws.Range("A189").Value = "Somalia"
$ws.Range("C189").Value = 1

$ws.Range("A190").Value = "Mauritania"
$ws.Range("C190").Value = 0

$ws.Range("A191").Value = "Republica del Chad"
$ws.Range("E191").Value = 3
$ws.Range("H191").Value = 0

$ws.Range("A192").Value = "Gambia"
$ws.Range("D192").Value = 0
$ws.Range("H192").Value = 1
